$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "K" (strikeout) values replacing the old "Strike#" values in column G,
# for rows 2-27.
$kValues = @(1,2,3,2,5,5,4,7,2,5,2,5,3,3,2,3,1,4,2,2,6,4,2,6,3,1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
